$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.359208748987044
$ws.Range("C2").Value = 0.3423610379078639
$ws.Range("E2").Value = 0.2727917404645233
$ws.Range("F2").Value = 1.486442873978078
$ws.Range("G2").Value = 0.3616040099458004
$ws.Range("H2").Value = 0.5373195866125116
$ws.Range("I2").Value = 0.5360708483043499
$ws.Range("J2").Value = 0.0315407555314664
$ws.Range("L2").Value = 0.5956268295430647
$ws.Range("O2").Value = 1.736698922457464

$ws.Range("B3").Value = 1.214844311253898
$ws.Range("C3").Value = 0.3316009879988258
$ws.Range("E3").Value = 0.2738473555045537
$ws.Range("F3").Value = 1.492251695326914
$ws.Range("G3").Value = 0.3662242352030844
$ws.Range("H3").Value = 0.5436084768193012
$ws.Range("I3").Value = 0.5464334666670112
$ws.Range("J3").Value = 0.02914222763023133
$ws.Range("L3").Value = 0.5784017149985914
$ws.Range("O3").Value = 1.75940046871257

$ws.Range("B4").Value = 1.125950509309348
$ws.Range("C4").Value = 0.3250177688963731
$ws.Range("E4").Value = 0.2746042062860958
$ws.Range("F4").Value = 1.496754891604297
$ws.Range("G4").Value = 0.369429427066656
$ws.Range("H4").Value = 0.5477771251002537
$ws.Range("I4").Value = 0.5532015393280236
$ws.Range("J4").Value = 0.02766495176064154
$ws.Range("L4").Value = 0.5679901586724156
$ws.Range("O4").Value = 1.774754585876451

$ws.Range("B5").Value = 1.089664578146994
$ws.Range("C5").Value = 0.322341253783307
$ws.Range("E5").Value = 0.2749400042547876
$ws.Range("F5").Value = 1.498825471451923
$ws.Range("G5").Value = 0.370827967521052
$ws.Range("H5").Value = 0.5495531195105272
$ws.Range("I5").Value = 0.5560613817296858
$ws.Range("J5").Value = 0.02706184680994284
$ws.Range("L5").Value = 0.5637891684078511
$ws.Range("O5").Value = 1.781366861283161

$ws.Range("B6").Value = 1.083635724825626
$ws.Range("C6").Value = 0.32189720360401
$ws.Range("E6").Value = 0.2749974178738945
$ws.Range("F6").Value = 1.499183512303127
$ws.Range("G6").Value = 0.3710657680755602
$ws.Range("H6").Value = 0.5498526862504178
$ws.Range("I6").Value = 0.5565424006362765
$ws.Range("J6").Value = 0.02696163643323501
$ws.Range("L6").Value = 0.5630941332058228
$ws.Range("O6").Value = 1.782486270995705

$ws.Range("B7").Value = 1.125461387821417
$ws.Range("C7").Value = 0.3249816469846252
$ws.Range("E7").Value = 0.274608624078212
$ws.Range("F7").Value = 1.496781862694256
$ws.Range("G7").Value = 0.3694479144553213
$ws.Range("H7").Value = 0.5478007640854301
$ws.Range("I7").Value = 0.5532396961720103
$ws.Range("J7").Value = 0.02765682248132606
$ws.Range("L7").Value = 0.5679333328826601
$ws.Range("O7").Value = 1.774842323256777

$ws.Range("B8").Value = 1.309486093623548
$ws.Range("C8").Value = 0.3386462671391826
$ws.Range("E8").Value = 0.2731331809673065
$ws.Range("F8").Value = 1.48825136243218
$ws.Range("G8").Value = 0.3631204752883548
$ws.Range("H8").Value = 0.5394242141552539
$ws.Range("I8").Value = 0.5395596595839631
$ws.Range("J8").Value = 0.03071471851670537
$ws.Range("L8").Value = 0.589653640869912
$ws.Range("O8").Value = 1.744232354960033

$ws.Range("B9").Value = 1.668243992941143
$ws.Range("C9").Value = 0.3656176118656163
$ws.Range("E9").Value = 0.2711005999179399
$ws.Range("F9").Value = 1.478956328911593
$ws.Range("G9").Value = 0.3536454246680591
$ws.Range("H9").Value = 0.5254366389490954
$ws.Range("I9").Value = 0.5159555261824629
$ws.Range("J9").Value = 0.03667324172383957
$ws.Range("L9").Value = 0.6335400151493786
$ws.Range("O9").Value = 1.695460023156613

$ws.Range("B10").Value = 1.930422451630761
$ws.Range("C10").Value = 0.3855270885087805
$ws.Range("E10").Value = 0.2701298757042601
$ws.Range("F10").Value = 1.476663801580401
$ws.Range("G10").Value = 0.348486522036886
$ws.Range("H10").Value = 0.5166479981009502
$ws.Range("I10").Value = 0.5005854896698434
$ws.Range("J10").Value = 0.04102588416834863
$ws.Range("L10").Value = 0.6665561920653431
$ws.Range("O10").Value = 1.666521018797468

$ws.Range("B11").Value = 2.049368107174132
$ws.Range("C11").Value = 0.3946021558028008
$ws.Range("E11").Value = 0.2698013108133281
$ws.Range("F11").Value = 1.47660713590804
$ws.Range("G11").Value = 0.3465340533700996
$ws.Range("H11").Value = 0.5129732044565571
$ws.Range("I11").Value = 0.4940230209965417
$ws.Range("J11").Value = 0.04300018848158516
$ws.Range("L11").Value = 0.6817409171869144
$ws.Range("O11").Value = 1.654860055898666

$ws.Range("B12").Value = 2.094361248421365
$ws.Range("C12").Value = 0.3980409994819354
$ws.Range("E12").Value = 0.2696931054737597
$ws.Range("F12").Value = 1.476727552289475
$ws.Range("G12").Value = 0.3458516560455891
$ws.Range("H12").Value = 0.5116281606430917
$ws.Range("I12").Value = 0.4915998839365869
$ws.Range("J12").Value = 0.04374694067493579
$ws.Range("L12").Value = 0.6875144454491817
$ws.Range("O12").Value = 1.650661146430664

$ws.Range("B13").Value = 2.084673395342747
$ws.Range("C13").Value = 0.3973002848022134
$ws.Range("E13").Value = 0.2697156888570937
$ws.Range("F13").Value = 1.476695307673381
$ws.Range("G13").Value = 0.3459960851706398
$ws.Range("H13").Value = 0.5119157701634691
$ws.Range("I13").Value = 0.492118992843773
$ws.Range("J13").Value = 0.04358615388419196
$ws.Range("L13").Value = 0.686269978304864
$ws.Range("O13").Value = 1.651555802410257

$ws.Range("B14").Value = 2.05307071842617
$ws.Range("C14").Value = 0.3948850268795923
$ws.Range("E14").Value = 0.2697920839355064
$ws.Range("F14").Value = 1.476614198940084
$ws.Range("G14").Value = 0.3464767691146093
$ws.Range("H14").Value = 0.5128616142919142
$ws.Range("I14").Value = 0.493822425935976
$ws.Range("J14").Value = 0.04306164203246254
$ws.Range("L14").Value = 0.6822154420780748
$ws.Range("O14").Value = 1.654510259575403

$ws.Range("B15").Value = 2.033706700080529
$ws.Range("C15").Value = 0.3934059040863929
$ws.Range("E15").Value = 0.2698409886944049
$ws.Range("F15").Value = 1.47658299540619
$ws.Range("G15").Value = 0.3467786273218181
$ws.Range("H15").Value = 0.5134470308294397
$ws.Range("I15").Value = 0.4948738978428171
$ws.Range("J15").Value = 0.04274024809599553
$ws.Range("L15").Value = 0.6797349575109877
$ws.Range("O15").Value = 1.656348208243884

$ws.Range("B16").Value = 1.92264237368903
$ws.Range("C16").Value = 0.3849343512279688
$ws.Range("E16").Value = 0.2701536192930156
$ws.Range("F16").Value = 1.476687355007286
$ws.Range("G16").Value = 0.3486220774057358
$ws.Range("H16").Value = 0.5168946619975898
$ws.Range("I16").Value = 0.5010230187719547
$ws.Range("J16").Value = 0.04089673911762759
$ws.Range("L16").Value = 0.6655671309479203
$ws.Range("O16").Value = 1.667313399727107

$ws.Range("B17").Value = 1.854423820045383
$ws.Range("C17").Value = 0.3797417635168756
$ws.Range("E17").Value = 0.2703743291613847
$ws.Range("F17").Value = 1.477004011860231
$ws.Range("G17").Value = 0.3498541563011699
$ws.Range("H17").Value = 0.519092483927821
$ws.Range("I17").Value = 0.5049054191785896
$ws.Range("J17").Value = 0.03976430210272497
$ws.Range("L17").Value = 0.6569177344587445
$ws.Range("O17").Value = 1.674425736105135

$ws.Range("B18").Value = 1.815156303365882
$ws.Range("C18").Value = 0.3767568499780793
$ws.Range("E18").Value = 0.2705119160844554
$ws.Range("F18").Value = 1.477278968592202
$ws.Range("G18").Value = 0.3505999177489514
$ws.Range("H18").Value = 0.5203870317334776
$ws.Range("I18").Value = 0.507178883533955
$ws.Range("J18").Value = 0.03911241745595362
$ws.Range("L18").Value = 0.6519584316796738
$ws.Range("O18").Value = 1.678658045909359

$ws.Range("B19").Value = 1.801855931280159
$ws.Range("C19").Value = 0.3757465150053747
$ws.Range("E19").Value = 0.2705603292995455
$ws.Range("F19").Value = 1.477388005042698
$ws.Range("G19").Value = 0.3508587843313222
$ws.Range("H19").Value = 0.5208305664678647
$ws.Range("I19").Value = 0.5079555746816755
$ws.Range("J19").Value = 0.03889160982348017
$ws.Range("L19").Value = 0.6502819886115105
$ws.Range("O19").Value = 1.68011531346616

$ws.Range("B20").Value = 1.861688928352066
$ws.Range("C20").Value = 0.3802943469128763
$ws.Range("E20").Value = 0.2703497332443519
$ws.Range("F20").Value = 1.476960695826548
$ws.Range("G20").Value = 0.3497191576073035
$ws.Range("H20").Value = 0.5188553734792833
$ws.Range("I20").Value = 0.5044879474422217
$ws.Range("J20").Value = 0.03988490787465793
$ws.Range("L20").Value = 0.6578368656801103
$ws.Range("O20").Value = 1.673653968261632

$ws.Range("B21").Value = 2.062354543377012
$ws.Range("C21").Value = 0.395594386263781
$ws.Range("E21").Value = 0.2697692050756615
$ws.Range("F21").Value = 1.476634171663022
$ws.Range("G21").Value = 0.3463340326879774
$ws.Range("H21").Value = 0.5125825340291357
$ws.Range("I21").Value = 0.4933204042632298
$ws.Range("J21").Value = 0.0432157279901233
$ws.Range("L21").Value = 0.6834057256689618
$ws.Range("O21").Value = 1.653636573743825

$ws.Range("B22").Value = 2.19321420577154
$ws.Range("C22").Value = 0.4056071908425452
$ws.Range("E22").Value = 0.2694842927239662
$ws.Range("F22").Value = 1.477247746953708
$ws.Range("G22").Value = 0.3444537820860489
$ws.Range("H22").Value = 0.5087540500689229
$ws.Range("I22").Value = 0.4863828102103227
$ws.Range("J22").Value = 0.04538749437723055
$ws.Range("L22").Value = 0.7002526779357936
$ws.Range("O22").Value = 1.641818275264256

$ws.Range("B23").Value = 2.123399142904532
$ws.Range("C23").Value = 0.4002620448324024
$ws.Range("E23").Value = 0.2696277217544605
$ws.Range("F23").Value = 1.476844583383624
$ws.Range("G23").Value = 0.3454268343020459
$ws.Range("H23").Value = 0.5107725564367414
$ws.Range("I23").Value = 0.4900524418576389
$ws.Range("J23").Value = 0.04422886625487621
$ws.Range("L23").Value = 0.6912488085407915
$ws.Range("O23").Value = 1.648010037359896

$ws.Range("B24").Value = 1.858404521438842
$ws.Range("C24").Value = 0.380044522771982
$ws.Range("E24").Value = 0.270360819733007
$ws.Range("F24").Value = 1.476979989614108
$ws.Range("G24").Value = 0.3497800739668193
$ws.Range("H24").Value = 0.5189624745816772
$ws.Range("I24").Value = 0.5046765574045455
$ws.Range("J24").Value = 0.03983038459162458
$ws.Range("L24").Value = 0.6574212848414049
$ws.Range("O24").Value = 1.674002438079

$ws.Range("B25").Value = 1.571428484211879
$ws.Range("C25").Value = 0.3583037725921656
$ws.Range("E25").Value = 0.2715585368850704
$ws.Range("F25").Value = 1.480674572833436
$ws.Range("G25").Value = 0.3558931843983615
$ws.Range("H25").Value = 0.5289594254674839
$ws.Range("I25").Value = 0.5219953180875958
$ws.Range("J25").Value = 0.03506557083183282
$ws.Range("L25").Value = 0.6215307638642003
$ws.Range("O25").Value = 1.707445969192307
